$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 519.5399500817239
$ws.Range("C2").Value = 981.181658104533
$ws.Range("D2").Value = 1329.718121219282
$ws.Range("E2").Value = 1482.659367292297
$ws.Range("B3").Value = 542.3516677527452
$ws.Range("C3").Value = 1001.504143663945
$ws.Range("D3").Value = 1340.620711963814
$ws.Range("E3").Value = 1489.660156732298
$ws.Range("B4").Value = 563.8605507489592
$ws.Range("C4").Value = 1039.496240544647
$ws.Range("D4").Value = 1385.142876463339
$ws.Range("E4").Value = 1536.745462111212
$ws.Range("B5").Value = 662.9786698217429
$ws.Range("C5").Value = 1131.549665732404
$ws.Range("D5").Value = 1436.697173459555
$ws.Range("E5").Value = 1589.172341735097
$ws.Range("B6").Value = 647.8814483192817
$ws.Range("C6").Value = 1122.600416929733
$ws.Range("D6").Value = 1432.230304965013
$ws.Range("E6").Value = 1586.12702903763
$ws.Range("B7").Value = 713.9659208237927
$ws.Range("C7").Value = 1176.796944956915
$ws.Range("D7").Value = 1471.982866696806
$ws.Range("E7").Value = 1595.448448823696
$ws.Range("B8").Value = 285.6775235257514
$ws.Range("C8").Value = 722.3126791818364
$ws.Range("D8").Value = 1114.811042334584
$ws.Range("E8").Value = 1433.203986887945
$ws.Range("B9").Value = 606.6176074155284
$ws.Range("C9").Value = 1065.463844133058
$ws.Range("D9").Value = 1396.439384704664
$ws.Range("E9").Value = 1546.589363101275
$ws.Range("B10").Value = 683.4259752372587
$ws.Range("C10").Value = 1157.648520803099
$ws.Range("D10").Value = 1464.053630130737
$ws.Range("E10").Value = 1589.572658536279
$ws.Range("B11").Value = 686.6652809504774
$ws.Range("C11").Value = 1160.705684869071
$ws.Range("D11").Value = 1466.482912877761
$ws.Range("E11").Value = 1591.98943187592
$ws.Range("B12").Value = 700.8260013260174
$ws.Range("C12").Value = 1172.403438093525
$ws.Range("D12").Value = 1471.719828093797
$ws.Range("E12").Value = 1594.044308298733
$ws.Range("B13").Value = 681.6785216095993
$ws.Range("C13").Value = 1155.606495902057
$ws.Range("D13").Value = 1462.307684232907
$ws.Range("E13").Value = 1587.962461674932
